# Fruta / hortaliza, semanal
# Insert a new row 8 (pushing the existing rows 8-13 down to 9-14) and
# populate it with the new weekly price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 8..13 down to 9..14, creating a blank row 8.
$ws.Rows.Item(8).Insert()

# Fill the new row 8 with the latest weekly data.
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = 'Vega Modelo de Temuco'
$ws.Range("C8").Value = 'La Araucanía'
$ws.Range("D8").Value = 44868
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 'Fruta'
$ws.Range("G8").Value = 100104
$ws.Range("H8").Value = 'Frutos de pepita'
$ws.Range("I8").Value = 100104004
$ws.Range("J8").Value = 'Níspero'
$ws.Range("K8").Value = 'Californiana(o)'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("Q8").Value = '$/bandeja 5 kilos'
$ws.Range("R8").Value = 'Provincia de Quillota'
$ws.Range("S8").Value = 2800
$ws.Range("T8").Value = 5
